# Weekly update: a new price record (row 8) is inserted for
# "Vega Monumental Concepción - Poroto granado", pushing all the
# subsequent records down by one row (old row 8 -> new row 9, ...,
# old row 25 -> new row 26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 8; Excel shifts rows 8:25 down to 9:26,
# carrying their existing values/formats with them.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with this week's record.
$ws.Cells.Item(8, 1).Value = 11
$ws.Cells.Item(8, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(8, 3).Value = "Bíobío"
$ws.Cells.Item(8, 4).Value = 44574
$ws.Cells.Item(8, 5).Value = 8
$ws.Cells.Item(8, 6).Value = 100112030
$ws.Cells.Item(8, 7).Value = "Poroto granado"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 100
$ws.Cells.Item(8, 11).Value = 30000
$ws.Cells.Item(8, 12).Value = 32000
$ws.Cells.Item(8, 13).Value = 31000
$ws.Cells.Item(8, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(8, 15).Value = "Región Metropolitana"
$ws.Cells.Item(8, 16).Value = 1240
$ws.Cells.Item(8, 17).Value = 25
$ws.Cells.Item(8, 18).Value = "Hortaliza"
